# Auto commit at 2025-09-28  8:36:10.73
# Append one new day (2025-09-27, serial 45927) of charging-station data
# as two new rows (54 and 55) at the bottom of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: date (stored as the Excel serial number so it keeps the
# existing yyyy-mm-dd number format already applied to the column).
$ws.Cells.Item(54, 1).Value = 45927
$ws.Cells.Item(55, 1).Value = 45927

# Column B: station name labels (reuse the existing text used elsewhere
# in the sheet).
$ws.Cells.Item(54, 2).Value = "四方坪站充电量(kw)"
$ws.Cells.Item(55, 2).Value = "高岭站充电量(kw)"

# Columns C:Z - hourly charging totals for each of the two stations.
$row54 = @(
    770.45399999999984,
    1074.317,
    458.66000000000008,
    710.23699999999997,
    391.58699999999993,
    683.07500000000005,
    434.56100000000009,
    213.12,
    301.85500000000002,
    229.68,
    247.38400000000001,
    69.656000000000006,
    864.14400000000023,
    952.14699999999993,
    465.35699999999997,
    386.322,
    390.822,
    172.47899999999998,
    149.07,
    53.460000000000008,
    309.45999999999998,
    75.671999999999997,
    238.70000000000002,
    47.25
)

$row55 = @(
    306.58899999999994,
    332.67,
    221.01899999999998,
    117.91999999999999,
    0,
    171.62,
    153.208,
    141.57999999999998,
    290.202,
    284.27300000000002,
    207.22900000000001,
    164.00799999999998,
    230.99400000000003,
    276.983,
    278.98599999999999,
    215.41200000000001,
    217.76799999999997,
    182.55,
    51.75,
    55.25,
    50.36,
    43.875,
    72.350000000000009,
    82.106999999999999
)

for ($i = 0; $i -lt $row54.Length; $i++) {
    $col = 3 + $i
    $ws.Cells.Item(54, $col).Value = $row54[$i]
    $ws.Cells.Item(55, $col).Value = $row55[$i]
}

# Match the author's final cursor position recorded in the saved view state.
[void]$ws.Range("E60").Select()
